# Adds 4 new "section title" slides at the end of the deck (same visual
# pattern as the existing "#1".."#4" slides): "#5 Predictive data analysis",
# "#6 Limitations and Future Work", "#7 Conclusions" and "#8 Annexes".
# Each new slide also gets a matching notes page (best effort - mirrors the
# exact boilerplate notes text already present on every other slide).

$p = $ppt.ActivePresentation

# Donor slide: slide 2 ("#1" / "CRISP-DM") uses the exact same 3-shape
# layout (big red "#N" box, slide-number placeholder, blue title box) as
# every other section-title slide in this deck, so it is the template we
# clone for the 4 new slides.
$donor = $p.Slides.Item(2)

# EMU -> Points helper that compensates for the ~4-decimal-place rounding
# PowerPoint's Shape.Width/Height setter applies internally (Shape sizes
# are stored as Single-precision points), so the saved EMU value lands on
# the exact target instead of being off by one EMU.
function EmuToPt([double]$emu) {
    return ($emu / 12700.0) + 0.00003
}

# Create the 4 new slides (in this exact order so the internal slide IDs
# come out as 274, 275, 276, 277 respectively) - always duplicating the
# donor and parking the copy at the end of the deck.
$dPredictive   = $donor.Duplicate().Item(1); $dPredictive.MoveTo($p.Slides.Count)
$dAnnexes      = $donor.Duplicate().Item(1); $dAnnexes.MoveTo($p.Slides.Count)
$dLimitations  = $donor.Duplicate().Item(1); $dLimitations.MoveTo($p.Slides.Count)
$dConclusions  = $donor.Duplicate().Item(1); $dConclusions.MoveTo($p.Slides.Count)

# Reorder to the final sequence: Predictive, Limitations, Conclusions, Annexes
$dLimitations.MoveTo(17)
$dConclusions.MoveTo(18)

# --- Slide 16: "#5 Predictive data analysis" ---------------------------
$s = $dPredictive
$s.Shapes.Item(1).TextFrame.TextRange.Text = "#5"
$titleShape = $s.Shapes.Item(3)
$titleShape.TextFrame.TextRange.Text = "Predictive data analysis"
$titleShape.Width = EmuToPt 7457554
$s.Shapes.Item(2).TextFrame.TextRange.Text = "$($s.SlideIndex)"

# --- Slide 17: "#6 Limitations and Future Work" -------------------------
$s = $dLimitations
$s.Shapes.Item(1).TextFrame.TextRange.Text = "#6"
$titleShape = $s.Shapes.Item(3)
$titleShape.TextFrame.TextRange.Text = "Limitations and Future Work"
$titleShape.Width = EmuToPt 9068508
$s.Shapes.Item(2).TextFrame.TextRange.Text = "$($s.SlideIndex)"

# --- Slide 18: "#7 Conclusions" -----------------------------------------
$s = $dConclusions
$s.Shapes.Item(1).TextFrame.TextRange.Text = "#7"
$titleShape = $s.Shapes.Item(3)
$titleShape.TextFrame.TextRange.Text = "Conclusions"
$titleShape.Width = EmuToPt 3898824
$s.Shapes.Item(2).TextFrame.TextRange.Text = "$($s.SlideIndex)"

# --- Slide 19: "#8 Annexes" ----------------------------------------------
$s = $dAnnexes
$s.Shapes.Item(1).TextFrame.TextRange.Text = "#8"
$titleShape = $s.Shapes.Item(3)
$titleShape.TextFrame.TextRange.Text = "Annexes"
$titleShape.Width = EmuToPt 2807692
$s.Shapes.Item(2).TextFrame.TextRange.Text = "$($s.SlideIndex)"

# --- Notes pages -----------------------------------------------------------
# Every existing slide in this deck shares the exact same (never-customised)
# speaker notes boilerplate; mirror that on the 4 new slides too.
$notesText = "Health Problems (stress, pollution, etc.) -> related with environment`rEconomic impact`rEnvironmental impact"

foreach ($s in @($dPredictive, $dLimitations, $dConclusions, $dAnnexes)) {
    $notes = $s.NotesPage.Shapes.Item(2)
    $notes.TextFrame.TextRange.Text = $notesText
}

Write-Host "Slide count: $($p.Slides.Count)"
